$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 329
$ws1.Range("F4").Value = 1441
$ws1.Range("F5").Value = 8630
$ws1.Range("F6").Value = 79
$ws1.Range("F9").Value = 271
$ws1.Range("F11").Value = 3480
$ws1.Range("F13").Value = 353
$ws1.Range("F15").Value = 1067
$ws1.Range("F20").Value = 2206
$ws1.Range("F21").Value = 8

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 329
$ws4.Range("F4").Value = 1441
$ws4.Range("F5").Value = 8630
$ws4.Range("F6").Value = 79
$ws4.Range("F11").Value = 3480
$ws4.Range("F13").Value = 353
$ws4.Range("F15").Value = 1067
$ws4.Range("F20").Value = 2206
$ws4.Range("F22").Value = 8
